$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "Complex Chaos"
$ws.Range("D8").Value = 1

$ws.Range("E10").Select()
